# Fruta / hortaliza, semanal
# Insert a new weekly record by shifting rows 399..419 down into 400..420
# (row 399's old values move to row 400, ..., row 419's old values move to
# row 420), then fill row 399 with the values that sat in row 398 before
# the shift. Only columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg)
# actually vary row-to-row for this product/market block - every other
# column is constant across the block, so a straight shift of those six
# columns (plus appending a brand-new row 420) reproduces the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that were in rows 398..419 (in that order) for columns D,J,K,L,M,P
# before the edit - index 0 goes to row 399, index 1 to row 400, ...,
# index 21 (the old row 419 values) to the brand-new row 420.
$D = @(45012,44705,44448,44970,45015,44342,44186,44179,44172,44441,44727,44855,44391,44168,44875,44642,44454,44426,44918,44243,44217,45007)
$J = @(500,500,400,500,500,500,500,400,500,500,500,500,400,500,500,400,350,400,500,800,500,200)
$K = @(2500,3000,2500,2500,2500,1800,1800,1500,1500,2000,2500,2500,2000,1500,2000,3500,2500,2000,2500,1500,1800,3000)
$L = @(2500,3000,2500,2500,2500,1800,1800,1500,1500,2000,2500,2500,2000,1500,2000,3500,2500,2000,2500,1500,1800,3000)
$M = @(2500,3000,2500,2500,2500,1800,1800,1500,1500,2000,2500,2500,2000,1500,2000,3500,2500,2000,2500,1500,1800,3000)
$P = @(625,750,625,625,625,450,450,375,375,500,625,625,500,375,500,875,625,500,625,375,450,750)

$startRow = 399

for ($i = 0; $i -lt $D.Length; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 4).Value = $D[$i]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 10).Value = $J[$i]
    $ws.Cells.Item($r, 11).Value = $K[$i]
    $ws.Cells.Item($r, 12).Value = $L[$i]
    $ws.Cells.Item($r, 13).Value = $M[$i]
    $ws.Cells.Item($r, 16).Value = $P[$i]
}

# Row 420 is brand new - populate the columns that were untouched by the
# per-cell loop above (they are constant across the whole block, so copy
# them straight from row 419, which still holds the same static values).
$lastRow = $startRow + $D.Length - 1

$ws.Cells.Item($lastRow, 1).Value = $ws.Cells.Item($lastRow - 1, 1).Value2
$ws.Cells.Item($lastRow, 2).Value = $ws.Cells.Item($lastRow - 1, 2).Value2
$ws.Cells.Item($lastRow, 3).Value = $ws.Cells.Item($lastRow - 1, 3).Value2
$ws.Cells.Item($lastRow, 5).Value = $ws.Cells.Item($lastRow - 1, 5).Value2
$ws.Cells.Item($lastRow, 6).Value = $ws.Cells.Item($lastRow - 1, 6).Value2
$ws.Cells.Item($lastRow, 7).Value = $ws.Cells.Item($lastRow - 1, 7).Value2
$ws.Cells.Item($lastRow, 8).Value = $ws.Cells.Item($lastRow - 1, 8).Value2
$ws.Cells.Item($lastRow, 9).Value = $ws.Cells.Item($lastRow - 1, 9).Value2
$ws.Cells.Item($lastRow, 14).Value = $ws.Cells.Item($lastRow - 1, 14).Value2
$ws.Cells.Item($lastRow, 15).Value = $ws.Cells.Item($lastRow - 1, 15).Value2
$ws.Cells.Item($lastRow, 17).Value = $ws.Cells.Item($lastRow - 1, 17).Value2
$ws.Cells.Item($lastRow, 18).Value = $ws.Cells.Item($lastRow - 1, 18).Value2

Write-Output "Rows $startRow..$lastRow updated; dimension now $($ws.UsedRange.Address())"
